# The author re-saved this workbook from Excel. Comparing the OOXML before/after,
# the only meaningful content edit is in the header row: the label in H1 changed
# from "Creditos" to "creditos" (the old shared string is dropped and a new one
# is appended, which is exactly what happens when a cell's text is edited in
# place). The rest of the diff (revisionPtr GUID, workbook window geometry,
# scrolled topLeftCell, and the shared-string index shuffle) is Excel's own
# save-time bookkeeping / view-state churn, not a user edit to replicate.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content edit: header cell H1 "Creditos" -> "creditos"
$ws.Range("H1").Value = "creditos"

# View-state: the saved file also shows the selection moved to I7 (from I15)
# and the sheet scrolled right so column E is left-most visible. The selection
# is reproducible via COM; apply it for parity.
$ws.Range("I7").Select()
